# Product Burndown sheet - "Edit Tutor accept/decline final bug fix"
#
# Sprint row 22 (2013-12-09) actually had plan-change / accomplished /
# actual-hours data that was missing (C22/D22/E22 were blank), and the next
# sprint's planned hours (B23) were mis-entered as 8 instead of 12. Enter the
# corrected inputs; every other cell on the sheet (F:K for rows 22-33, plus
# the D34/E34 sprint averages, plus the burndown chart that reads off column
# J) is formula-driven and recalculates automatically from these four cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 13
$ws.Range("B23").Value = 12

# Match the author's final cursor position on the sheet.
$ws.Range("J26").Select()
